$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("model_4_2_24", 0.6402165563634941, -5.835797479179137, -0.08935279892930792, -2.234412954938043, 0.3981742858886719, 1.664963960647583, 0.487859845161438, 1.111031055450439),
    @("model_4_2_23", 0.642067167262961, -5.810251914785142, -0.07116339701319041, -2.213674993957021, 0.3961262106895447, 1.65874195098877, 0.4797138571739197, 1.103907465934753),
    @("model_4_2_22", 0.6455745112936551, -5.801301765222518, -0.01438509170510516, -2.175471984170538, 0.3922445774078369, 1.656561970710754, 0.4542860388755798, 1.090784549713135),
    @("model_4_2_21", 0.6501628585389463, -5.711535073065461, 0.01142570791889974, -2.125940142440476, 0.3871666491031647, 1.634697914123535, 0.4427268505096436, 1.07377016544342),
    @("model_4_2_20", 0.6518119711407308, -5.66871311856145, 0.0159222653978961, -2.107100581472376, 0.3853415548801422, 1.624268054962158, 0.4407130777835846, 1.067298769950867),
    @("model_4_2_19", 0.6526514280987448, -5.716588096234124, 0.06172610472735185, -2.096973698990188, 0.3844124972820282, 1.635928630828857, 0.4202001094818115, 1.063820123672485),
    @("model_4_2_17", 0.6561275872804313, -5.670199159532011, 0.09583786499475444, -2.058636966909144, 0.3805654346942902, 1.624629855155945, 0.4049233794212341, 1.05065131187439),
    @("model_4_2_11", 0.6561759489033672, -5.681032381076133, 0.1094878240627549, -2.054322693729076, 0.3805119097232819, 1.627268433570862, 0.3988103270530701, 1.049169421195984),
    @("model_4_2_18", 0.6564048685590793, -5.645205317256752, 0.08436641545759427, -2.056291714512841, 0.3802585601806641, 1.618542194366455, 0.4100607931613922, 1.049845695495605),
    @("model_4_2_10", 0.6577669475779759, -5.682806371858902, 0.1408505612465881, -2.035760212857919, 0.3787511885166168, 1.627700448036194, 0.3847647309303284, 1.042793035507202),
    @("model_4_2_16", 0.6577733684925429, -5.656114074839001, 0.116517287306139, -2.040662241512929, 0.3787440359592438, 1.621199131011963, 0.395662248134613, 1.044476985931396),
    @("model_4_2_15", 0.6584930828626734, -5.64502335122896, 0.1227232427723687, -2.032686516407755, 0.3779475390911102, 1.618497967720032, 0.3928829431533813, 1.04173743724823),
    @("model_4_2_14", 0.6587273985817459, -5.632319424524577, 0.1193122809858661, -2.030009979952111, 0.3776881992816925, 1.615403652191162, 0.3944104909896851, 1.040817856788635),
    @("model_4_2_13", 0.6591758951755795, -5.609036931486633, 0.1135299995090239, -2.024820998622188, 0.3771918118000031, 1.609732866287231, 0.3970000743865967, 1.039035439491272),
    @("model_4_2_9", 0.6592590874564865, -5.617174287328808, 0.1275885968144411, -2.019255627091584, 0.3770997822284698, 1.611714839935303, 0.3907040059566498, 1.037123680114746),
    @("model_4_2_12", 0.6595737641292584, -5.60685984168388, 0.1195200679260849, -2.02032412076325, 0.3767515122890472, 1.60920262336731, 0.3943174481391907, 1.037490844726562),
    @("model_4_2_8", 0.6598349436076509, -5.626213187498289, 0.1455792234618556, -2.011603806331865, 0.3764624893665314, 1.613916516304016, 0.3826470375061035, 1.03449535369873),
    @("model_4_2_7", 0.6702017318967624, -5.500242178577265, 0.253084060473402, -1.898355812661447, 0.3649895191192627, 1.583234190940857, 0.3345016539096832, 0.9955942630767822),
    @("model_4_2_6", 0.677115221928447, -5.343893932930807, 0.2806650387691867, -1.822746524713773, 0.3573382794857025, 1.545153141021729, 0.3221496641635895, 0.9696221947669983),
    @("model_4_2_2", 0.6798095858891022, -4.750859394817433, 0.1083871093176918, -1.705824576852493, 0.3543564677238464, 1.400710463523865, 0.3993032574653625, 0.9294591546058655),
    @("model_4_2_4", 0.6812088425828139, -5.157150963158012, 0.2561068811118274, -1.767711852397233, 0.3528078496456146, 1.499669075012207, 0.3331479132175446, 0.9507176876068115),
    @("model_4_2_5", 0.6831110222607655, -5.223389458654476, 0.3182598937007072, -1.754448222787054, 0.350702702999115, 1.515802383422852, 0.3053130507469177, 0.9461616277694702),
    @("model_4_2_3", 0.688897010499768, -4.907401762844791, 0.2667207544175996, -1.667445455201377, 0.3442993462085724, 1.438838720321655, 0.3283945322036743, 0.9162758588790894),
    @("model_4_2_1", 0.7723003809953508, -1.726297400591362, 0.3822991686042894, -0.4023930227929842, 0.2519963681697845, 0.6640318036079407, 0.2766334712505341, 0.4817264080047607),
    @("model_4_2_0", 0.7826166503205511, -1.225931746026951, 0.4042675848874772, -0.2010835684546819, 0.240579292178154, 0.5421599149703979, 0.2667950689792633, 0.4125759601593018)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 1 + $j
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
